# "paises.xlsx" / sheet "Pais" daily COVID data refresh (25 Aug 22:50 -> 26 Aug 00:07).
#
# Column layout: A=Pais, B=Casos totales, C=Nuevos casos, D=Casos activos,
#                E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes
#
# Most rows simply get refreshed B:H statistics. A handful of countries
# (Ruanda, Siria, Togo, Guyana, Islas Malvinas) were re-ranked by the source
# and now occupy a row earlier than before, which pushes the countries that
# used to sit there down by one row - so those rows need both their A
# (country name) and B:H values rewritten to match the new day's snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Agosto de 2020 a las 00:07"

# Row 4
$ws.Range("B4").Value = 5949693
$ws.Range("C4").Value = 34063
$ws.Range("D4").Value = 3238686
$ws.Range("E4").Value = 2528780
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1113
$ws.Range("H4").Value = 182227

# Row 5
$ws.Range("B5").Value = 3669995
$ws.Range("C5").Value = 42778
$ws.Range("D5").Value = 2848395
$ws.Range("E5").Value = 705020
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 1129
$ws.Range("H5").Value = 116580

# Row 6
$ws.Range("B6").Value = 3231754
$ws.Range("C6").Value = 66873
$ws.Range("D6").Value = 2467252
$ws.Range("E6").Value = 704890
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 1066
$ws.Range("H6").Value = 59612

# Row 54
$ws.Range("B54").Value = 50076
$ws.Range("C54").Value = 357
$ws.Range("D54").Value = 46673
$ws.Range("E54").Value = 3217
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 1
$ws.Range("H54").Value = 186

# Row 79
$ws.Range("B79").Value = 17562
$ws.Range("C79").Value = 56
$ws.Range("D79").Value = 15908
$ws.Range("E79").Value = 1540
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 114

# Row 81
$ws.Range("B81").Value = 15589
$ws.Range("C81").Value = 203
$ws.Range("D81").Value = 10601
$ws.Range("E81").Value = 4416
$ws.Range("F81").Value = 0
$ws.Range("G81").Value = 9
$ws.Range("H81").Value = 572

# Row 93
$ws.Range("B93").Value = 9128
$ws.Range("C93").Value = 52
$ws.Range("D93").Value = 8040
$ws.Range("E93").Value = 1031
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 2
$ws.Range("H93").Value = 57

# Row 107
$ws.Range("B107").Value = 5423
$ws.Range("C107").Value = 4
$ws.Range("D107").Value = 3066
$ws.Range("E107").Value = 2187
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 170

# Row 119
$ws.Range("B119").Value = 3568
$ws.Range("C119").Value = 36
$ws.Range("D119").Value = 2673
$ws.Range("E119").Value = 858
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 37

# Row 120
$ws.Range("A120").Value = "Ruanda"
$ws.Range("B120").Value = 3537
$ws.Range("C120").Value = 231
$ws.Range("D120").Value = 1806
$ws.Range("E120").Value = 1716
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = 15

# Row 121
$ws.Range("A121").Value = "Mozambique"
$ws.Range("B121").Value = 3508
$ws.Range("C121").Value = 68
$ws.Range("D121").Value = 1809
$ws.Range("E121").Value = 1678
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 21

# Row 122
$ws.Range("A122").Value = "Eslovaquia"
$ws.Range("B122").Value = 3452
$ws.Range("C122").Value = 28
$ws.Range("D122").Value = 2167
$ws.Range("E122").Value = 1252
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 33

# Row 123
$ws.Range("A123").Value = "Tailandia"
$ws.Range("B123").Value = 3402
$ws.Range("C123").Value = 5
$ws.Range("D123").Value = 3229
$ws.Range("E123").Value = 115
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 58

# Row 124
$ws.Range("B124").Value = 3275
$ws.Range("C124").Value = 6
$ws.Range("D124").Value = 2443
$ws.Range("E124").Value = 737
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 2
$ws.Range("H124").Value = 95

# Row 134
$ws.Range("A134").Value = "Siria"
$ws.Range("B134").Value = 2365
$ws.Range("C134").Value = 72
$ws.Range("D134").Value = 533
$ws.Range("E134").Value = 1737
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 3
$ws.Range("H134").Value = 95

# Row 135
$ws.Range("A135").Value = "Estonia"
$ws.Range("B135").Value = 2294
$ws.Range("C135").Value = 19
$ws.Range("D135").Value = 2038
$ws.Range("E135").Value = 192
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 64

# Row 136
$ws.Range("B136").Value = 2283
$ws.Range("C136").Value = 61
$ws.Range("D136").Value = 977
$ws.Range("E136").Value = 1204
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 2
$ws.Range("H136").Value = 102

# Row 140
$ws.Range("B140").Value = 2001
$ws.Range("C140").Value = 4
$ws.Range("D140").Value = 1569
$ws.Range("E140").Value = 363
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 69

# Row 141
$ws.Range("B141").Value = 1924
$ws.Range("C141").Value = 8
$ws.Range("D141").Value = 1091
$ws.Range("E141").Value = 276
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 2
$ws.Range("H141").Value = 557

# Row 154
$ws.Range("A154").Value = "Togo"
$ws.Range("B154").Value = 1309
$ws.Range("C154").Value = 14
$ws.Range("D154").Value = 919
$ws.Range("E154").Value = 363
$ws.Range("F154").Value = 0
$ws.Range("G154").Value = 0
$ws.Range("H154").Value = 27

# Row 155
$ws.Range("A155").Value = "Liberia"
$ws.Range("B155").Value = 1295
$ws.Range("C155").Value = 5
$ws.Range("D155").Value = 821
$ws.Range("E155").Value = 392
$ws.Range("F155").Value = 0
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 82

# Row 158
$ws.Range("B158").Value = 1173
$ws.Range("C158").Value = 1
$ws.Range("D158").Value = 1084
$ws.Range("E158").Value = 20
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 69

# Row 159
$ws.Range("A159").Value = "Guyana"
$ws.Range("B159").Value = 1060
$ws.Range("C159").Value = 31
$ws.Range("D159").Value = 523
$ws.Range("E159").Value = 506
$ws.Range("F159").Value = 0
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 31

# Row 160
$ws.Range("A160").Value = "Principado de Andorra"
$ws.Range("B160").Value = 1060
$ws.Range("C160").Value = 0
$ws.Range("D160").Value = 877
$ws.Range("E160").Value = 130
$ws.Range("F160").Value = 0
$ws.Range("G160").Value = 0
$ws.Range("H160").Value = 53

# Row 161
$ws.Range("A161").Value = "Lesoto"
$ws.Range("B161").Value = 1049
$ws.Range("C161").Value = 34
$ws.Range("D161").Value = 472
$ws.Range("E161").Value = 547
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 0
$ws.Range("H161").Value = 30

# Row 214
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

# Row 215
$ws.Range("A215").Value = "Montserrat"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 12
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 1
